# Update the lattice-multiplication problem table: replace the 15
# "A x B" exercises (and their derived scaffold lines) with a new set of
# problems, in row-major order (5 rows x 3 columns), leaving the table
# structure, cell formatting (sz=32) and paragraph/run shape untouched.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# New "A x B" problems, in row-major (top-left to bottom-right) order.
$problems = @(
    @("14", "79"),
    @("39", "39"),
    @("75", "16"),
    @("45", "11"),
    @("52", "71"),
    @("44", "99"),
    @("16", "13"),
    @("53", "43"),
    @("82", "59"),
    @("13", "78"),
    @("78", "67"),
    @("99", "52"),
    @("76", "65"),
    @("27", "82"),
    @("47", "49")
)

# NOTE: PowerShell's "+" operator performs *numeric* addition when both
# operands look like numbers (even single-digit strings pulled from
# string indexing), so string pieces are combined with the "-f" format
# operator / -join instead of "+" to guarantee textual concatenation.
$nl = [char]11  # maps to <w:br/> when assigned via Range.Text
$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count

$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $pair = $problems[$i]
        $aStr = $pair[0]
        $bStr = $pair[1]

        $header = "{0} x {1}" -f $aStr, $bStr
        $multiplierLine = "  {0}    {1}" -f $bStr[0], $bStr[1]
        $dashLine = "  ----"
        $row1 = "{0}|    |" -f $aStr[0]
        $row2 = "{0}|    |" -f $aStr[1]

        $lines = @($header, $multiplierLine, $dashLine, $row1, $row2)
        $newText = $lines -join $nl

        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $newText

        $i++
    }
}
